# Commit: "Act greficos y tablas web pob"
#
# 1. Rename sheets: "Datos" -> "Data", "Ficha técnica" -> "Metadata"
# 2. Reverse the chronological order of the data table in the "Data" sheet
#    (years 2006-2019 were ascending top to bottom; now descending 2019-2006).
#    Column A (year) is text, so the rows are reordered with Copy() (rather
#    than re-typing the values) to avoid Excel re-interpreting the year text
#    as a number.
# 3. Rework the "Metadata" sheet: lowercase the field-name keys, add a new
#    "observaciones" (Sin observaciones) row, and append a final
#    attribution row ("Mirador DESCA - UMAD/FCS - INDDHH").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename worksheets
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsMeta = $wb.Worksheets.Item(2)

$wsData.Name = "Data"
$wsMeta.Name = "Metadata"

# ---------------------------------------------------------------------------
# 2. Reverse the year rows (2..15) in the Data sheet, column A through O
# ---------------------------------------------------------------------------
$firstRow = 2
$lastRow = 15
$stagingOffset = 200   # far enough away to not collide with real data

# Stage a copy of every row so the originals are preserved while we
# overwrite rows 2..15 in reverse order.
for ($i = $firstRow; $i -le $lastRow; $i++) {
    $stagedRow = $stagingOffset + $i
    $srcRange = "A${i}:O${i}"
    $dstRange = "A${stagedRow}:O${stagedRow}"
    $wsData.Range($srcRange).Copy($wsData.Range($dstRange))
}

# Copy back from staging into rows 2..15, reversing the row order.
for ($i = $firstRow; $i -le $lastRow; $i++) {
    $srcOldRow = ($firstRow + $lastRow) - $i
    $stagedRow = $stagingOffset + $srcOldRow
    $srcRange = "A${stagedRow}:O${stagedRow}"
    $dstRange = "A${i}:O${i}"
    $wsData.Range($srcRange).Copy($wsData.Range($dstRange))
}

# Remove the temporary staging rows.
$clearFirst = $stagingOffset + $firstRow
$clearLast = $stagingOffset + $lastRow
$wsData.Range("A${clearFirst}:O${clearLast}").Clear()

# ---------------------------------------------------------------------------
# 3. Rework the Metadata sheet contents
# ---------------------------------------------------------------------------
# Row 1 (blank key / single-space value) is left untouched - it is unchanged
# by the commit.

$wsMeta.Cells.Item(2, 1).Value = "nomindicador"
$wsMeta.Cells.Item(2, 2).Value = "Porcentaje de personas de 4 a 17 años que no asisten a centros educativos"

$wsMeta.Cells.Item(3, 1).Value = "derecho"
$wsMeta.Cells.Item(3, 2).Value = "Educación"

$wsMeta.Cells.Item(4, 1).Value = "conindicador"
$wsMeta.Cells.Item(4, 2).Value = "No asistencia a centros educativos (4 a 17 años)"

$wsMeta.Cells.Item(5, 1).Value = "tipoind"
$wsMeta.Cells.Item(5, 2).Value = "Resultados"

$wsMeta.Cells.Item(6, 1).Value = "definicion"
$wsMeta.Cells.Item(6, 2).Value = "El indicador mide el porcentaje de personas de 4 a 17 años que no asisten a centros educativos."

$wsMeta.Cells.Item(7, 1).Value = "calculo"
$wsMeta.Cells.Item(7, 2).Value = "Para cada año calcular: (Cantidad de niños de 4 a 17 años que no asisten a centros educativos / Cantidad de población de 4 a 17 años)*100"

$wsMeta.Cells.Item(8, 1).Value = "observaciones"
$wsMeta.Cells.Item(8, 2).Value = "Sin observaciones"

$wsMeta.Cells.Item(9, 1).Value = "cita"
$wsMeta.Cells.Item(9, 2).Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"

$wsMeta.Cells.Item(10, 1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$wsMeta.Cells.Item(10, 2).Value = " "
